$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("books")

# Swap the header values in A1 and B1 ("author" <-> "title")
$a1 = $ws.Range("A1").Value2
$b1 = $ws.Range("B1").Value2
$ws.Range("A1").Value2 = $b1
$ws.Range("B1").Value2 = $a1

# Move the active selection to B1
$ws.Range("B1").Select()
